# Add two new columns, I ("I0") and J ("IF"), to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells - copy the formatting of the existing header cell (H1, style index 1:
# bold font, thin border, centered/top aligned) onto the new header cells, then set
# their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:J47 (column I = "I0", column J = "IF")
$IVals = @(8,6,8,7,7,7,8,7,9,8,7,7,6,8,8,8,8,8,8,9,5,7,7,8,6,6,6,7,7,8,5,7,7,9,8,8,8,6,6,8,8,8,7,9,8,5)
$JVals = @(8,6,8,7,7,7,8,7,9,8,7,7,6,8,8,8,8,8,8,9,6,7,8,8,6,6,6,7,8,8,6,7,7,9,8,8,8,6,6,8,9,8,7,9,8,6)

for ($i = 0; $i -lt $IVals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $IVals[$i]
    $ws.Cells.Item($row, 10).Value = $JVals[$i]
}
